$d = $word.ActiveDocument

function Split-WithBreak($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# --- Paragrafo "Programa" ---
# Divide o texto corrido em 5 blocos (introducao + itens 1-4), inserindo
# quebras de linha manuais (<w:br/>) entre eles, mantendo tudo no mesmo run.
Split-WithBreak "individual.1. Elaboração" "individual.^l1. Elaboração"
Split-WithBreak "digital2. Desenvolvimento" "digital^l2. Desenvolvimento"
Split-WithBreak "protótipos.3. O aluno" "protótipos.^l3. O aluno"
Split-WithBreak "IV).4. Apresentação" "IV).^l4. Apresentação"

# --- Paragrafo "Avaliacao" / "Metodo:" ---
Split-WithBreak "outros.Os alunos serão divididos" "outros.^lOs alunos serão divididos"
Split-WithBreak "profissão.Cada grupo deverá buscar" "profissão.^lCada grupo deverá buscar"
Split-WithBreak "projeto.As aulas ocorrerão" "projeto.^lAs aulas ocorrerão"

# --- Paragrafo "Avaliacao" / "Criterio:" ---
Split-WithBreak "outros.O detalhamento dos pesos" "outros.^lO detalhamento dos pesos"
